$d = $word.ActiveDocument

# Locate the run ".: ${" that precedes the "protocol" merge field in the
# "Αρ. Πρωτ.: ${protocol}" line, and remember its start/end.
$rng = $d.Content
$found = $rng.Find.Execute(".: `${", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "target text '.: `${' not found"
}

$insertAt = $rng.Start + 3   # right after ".: " (3 chars), before "${"

# Insert the new protocol-prefix text "Φ.15.1/" at that point.
$insPoint = $d.Range($insertAt, $insertAt)
$insPoint.InsertAfter("Φ.15.1/")

# Force the newly inserted text into its own run (matching the formatting
# of the surrounding text) by round-tripping Bold, which materializes the
# full run properties (rFonts ascii/hAnsi/cs, sz, szCs) on the new range.
$newRng = $d.Range($insertAt, $insertAt + 7)
$newRng.Bold = 1
$newRng.Bold = 0

Write-Output "done"
